$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = "2024-06-15 21:12:05"
$ws.Range("D48").Value = 200
$ws.Range("E48").Value = 13

# Row 49
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 2
$ws.Range("C49").Value = "2024-06-15 21:12:05"
$ws.Range("D49").Value = 200
$ws.Range("E49").Value = 0
